$wb = $excel.ActiveWorkbook

# Sheet ALC, row 13
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 9000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 9000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 9000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -9338

# Sheet ALC, row 16
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 10000
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4168622.5
$ws.Range("J17").Value = 4168622.5
$ws.Range("L17").Value = 12505867.5
$ws.Range("N17").Value = -12506203.5

# Sheet ALC, row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1589.2727
$ws.Range("I28").Value = 1221.4706
$ws.Range("J28").Value = 2839.8
$ws.Range("K28").Value = 1221.4706
$ws.Range("L28").Value = 2839.8
$ws.Range("M28").Value = -736.4706000000001
$ws.Range("N28").Value = -3809.8

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2729.8333
$ws.Range("I62").Value = 3015.8
$ws.Range("J62").Value = 1300
$ws.Range("K62").Value = 3015.8
$ws.Range("L62").Value = 1300
$ws.Range("M62").Value = -2391.8
$ws.Range("N62").Value = -2548

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2729.8333
$ws.Range("I65").Value = 3015.8
$ws.Range("J65").Value = 1300
$ws.Range("K65").Value = 15079
$ws.Range("L65").Value = 6500
$ws.Range("M65").Value = -11959
$ws.Range("N65").Value = -12740

# Sheet ALC, row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1310.15
$ws.Range("I80").Value = 1347
$ws.Range("J80").Value = 1273.3
$ws.Range("K80").Value = 4041
$ws.Range("L80").Value = 3819.9
$ws.Range("M80").Value = -3043
$ws.Range("N80").Value = -5815.9

# Sheet ALC, row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1310.15
$ws.Range("I83").Value = 1347
$ws.Range("J83").Value = 1273.3
$ws.Range("K83").Value = 12123
$ws.Range("L83").Value = 11459.7
$ws.Range("M83").Value = -7131
$ws.Range("N83").Value = -21443.7

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2343.739
$ws.Range("I113").Value = 2447.0588
$ws.Range("J113").Value = 2051
$ws.Range("K113").Value = 2447.0588
$ws.Range("L113").Value = 2051
$ws.Range("M113").Value = 806.9412000000002
$ws.Range("N113").Value = -8559

# Sheet ALC, row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2095.3333
$ws.Range("I125").Value = 2399.4546
$ws.Range("J125").Value = 1259
$ws.Range("K125").Value = 21595.0914
$ws.Range("L125").Value = 11331
$ws.Range("M125").Value = -19135.0914
$ws.Range("N125").Value = -16251

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2325.58
$ws.Range("I138").Value = 1156.5
$ws.Range("J138").Value = 3104.9666
$ws.Range("K138").Value = 3469.5
$ws.Range("L138").Value = 9314.899800000001
$ws.Range("M138").Value = 1670.5
$ws.Range("N138").Value = -19594.8998

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1216.7742
$ws.Range("I141").Value = 1030.3704
$ws.Range("J141").Value = 2475
$ws.Range("K141").Value = 3091.1112
$ws.Range("L141").Value = 7425
$ws.Range("M141").Value = 2088.8888
$ws.Range("N141").Value = -17785

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1923
$ws.Range("I102").Value = 1075
$ws.Range("J102").Value = 2262.2
$ws.Range("K102").Value = 1075
$ws.Range("L102").Value = 2262.2
$ws.Range("M102").Value = 547
$ws.Range("N102").Value = -5506.2

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6413588.5
$ws.Range("I31").Value = 3412.6938
$ws.Range("K31").Value = 3412.6938
$ws.Range("M31").Value = -3117.6938

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6413588.5
$ws.Range("I34").Value = 3412.6938
$ws.Range("K34").Value = 3412.6938
$ws.Range("M34").Value = -3210.6938

# Sheet CRP, row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1485.7693
$ws.Range("I105").Value = 1694.1111
$ws.Range("J105").Value = 1017
$ws.Range("K105").Value = 1694.1111
$ws.Range("L105").Value = 1017
$ws.Range("M105").Value = 52.88889999999992
$ws.Range("N105").Value = -4511

# Sheet GSM, row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 114.4
$ws.Range("I2").Value = 123.333336
$ws.Range("J2").Value = 101
$ws.Range("K2").Value = 123.333336
$ws.Range("L2").Value = 101
$ws.Range("M2").Value = -10.333336
$ws.Range("N2").Value = -327

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2990.4285
$ws.Range("I97").Value = 2982.2
$ws.Range("K97").Value = 2982.2
$ws.Range("M97").Value = -2486.2

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3510834
$ws.Range("I122").Value = 4763311.5
$ws.Range("K122").Value = 14289934.5
$ws.Range("M122").Value = -14287484.5

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3344.2
$ws.Range("I132").Value = 2126.4583
$ws.Range("J132").Value = 6001.091
$ws.Range("K132").Value = 6379.374899999999
$ws.Range("L132").Value = 18003.273
$ws.Range("M132").Value = -3849.374899999999
$ws.Range("N132").Value = -23063.273

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11625
$ws.Range("I40").Value = 14750
$ws.Range("J40").Value = 8500
$ws.Range("K40").Value = 14750
$ws.Range("L40").Value = 8500
$ws.Range("M40").Value = -14614
$ws.Range("N40").Value = -8772

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2751.182
$ws.Range("I82").Value = 2200
$ws.Range("J82").Value = 2957.875
$ws.Range("K82").Value = 2200
$ws.Range("L82").Value = 2957.875
$ws.Range("M82").Value = -1839
$ws.Range("N82").Value = -3679.875

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2751.182
$ws.Range("I85").Value = 2200
$ws.Range("J85").Value = 2957.875
$ws.Range("K85").Value = 2200
$ws.Range("L85").Value = 2957.875
$ws.Range("M85").Value = -952
$ws.Range("N85").Value = -5453.875

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1148.3636
$ws.Range("I93").Value = 913.3333
$ws.Range("J93").Value = 1430.4
$ws.Range("K93").Value = 913.3333
$ws.Range("L93").Value = 1430.4
$ws.Range("M93").Value = 334.6667
$ws.Range("N93").Value = -3926.4

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2067.1428
$ws.Range("I100").Value = 2049.2222
$ws.Range("J100").Value = 2099.4
$ws.Range("K100").Value = 2049.2222
$ws.Range("L100").Value = 2099.4
$ws.Range("M100").Value = -1508.2222
$ws.Range("N100").Value = -3181.4

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4486.9414
$ws.Range("I122").Value = 4134.24
$ws.Range("K122").Value = 12402.72
$ws.Range("M122").Value = -9952.719999999999

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10920.033
$ws.Range("I132").Value = 7979.8667
$ws.Range("K132").Value = 23939.6001
$ws.Range("M132").Value = -21409.6001

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10874105
$ws.Range("I136").Value = 13159784
$ws.Range("K136").Value = 39479352
$ws.Range("M136").Value = -39476802

# Sheet WVR, row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 59919.332
$ws.Range("J46").Value = 59919.332
$ws.Range("L46").Value = 59919.332
$ws.Range("N46").Value = -60381.332

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 707.2857
$ws.Range("I81").Value = 349.5
$ws.Range("J81").Value = 850.4
$ws.Range("K81").Value = 699
$ws.Range("L81").Value = 1700.8
$ws.Range("M81").Value = 362
$ws.Range("N81").Value = -3822.8

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 707.2857
$ws.Range("I84").Value = 349.5
$ws.Range("J84").Value = 850.4
$ws.Range("K84").Value = 3495
$ws.Range("L84").Value = 8504
$ws.Range("M84").Value = 1809
$ws.Range("N84").Value = -19112

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3079.2222
$ws.Range("I122").Value = 3326
$ws.Range("J122").Value = 1105
$ws.Range("K122").Value = 9978
$ws.Range("L122").Value = 3315
$ws.Range("M122").Value = -7528
$ws.Range("N122").Value = -8215

# Sheet WVR, row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 59919.332
$ws.Range("J134").Value = 59919.332
$ws.Range("L134").Value = 179757.996
$ws.Range("N134").Value = -184827.996

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1538.8
$ws.Range("I136").Value = 1538.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4616.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2066.4
$ws.Range("N136").ClearContents()
